$d = $word.ActiveDocument

$findText = " add date, hours(or time in: time out), and a drop down selector for category of hours"

$range = $d.Content
$range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $findText, 2)
